$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 125095.25
$ws.Range("I61").Value = 127
$ws.Range("K61").Value = 381
$ws.Range("M61").Value = -209
$ws.Range("H68").Value = 29317.6
$ws.Range("J68").Value = 29317.6
$ws.Range("L68").Value = 29317.6
$ws.Range("N68").Value = -30815.6
$ws.Range("H71").Value = 29317.6
$ws.Range("J71").Value = 29317.6
$ws.Range("L71").Value = 87952.79999999999
$ws.Range("N71").Value = -95440.79999999999
$ws.Range("H98").Value = 4215.706
$ws.Range("I98").Value = 5024.273
$ws.Range("J98").Value = 2733.3333
$ws.Range("K98").Value = 5024.273
$ws.Range("L98").Value = 2733.3333
$ws.Range("M98").Value = -3526.273
$ws.Range("N98").Value = -5729.3333
$ws.Range("H112").Value = 4465848.5
$ws.Range("J112").Value = 4546827.5
$ws.Range("L112").Value = 13640482.5
$ws.Range("N112").Value = -13642698.5
$ws.Range("H116").Value = 3625.45
$ws.Range("I116").Value = 2987.6155
$ws.Range("K116").Value = 2987.6155
$ws.Range("M116").Value = 454.3845000000001
$ws.Range("H122").Value = 4215.706
$ws.Range("I122").Value = 5024.273
$ws.Range("J122").Value = 2733.3333
$ws.Range("K122").Value = 15072.819
$ws.Range("L122").Value = 8199.999899999999
$ws.Range("M122").Value = -12622.819
$ws.Range("N122").Value = -13099.9999
$ws.Range("H132").Value = 3638395.8
$ws.Range("I132").Value = 4001840.2
$ws.Range("K132").Value = 12005520.6
$ws.Range("M132").Value = -12002990.6
$ws.Range("H138").Value = 2506.2642
$ws.Range("I138").Value = 1538.6129
$ws.Range("J138").Value = 3869.7727
$ws.Range("K138").Value = 4615.8387
$ws.Range("L138").Value = 11609.3181
$ws.Range("M138").Value = 524.1612999999998
$ws.Range("N138").Value = -21889.3181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8408.299000000001
$ws.Range("I32").Value = 6902.239
$ws.Range("J32").Value = 18498.9
$ws.Range("K32").Value = 6902.239
$ws.Range("L32").Value = 18498.9
$ws.Range("M32").Value = -6615.239
$ws.Range("N32").Value = -19072.9
$ws.Range("H46").Value = 6223.3335
$ws.Range("J46").Value = 6223.3335
$ws.Range("L46").Value = 6223.3335
$ws.Range("N46").Value = -6861.3335
$ws.Range("H63").Value = 2016.6666
$ws.Range("I63").Value = 2016.6666
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2016.6666
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1330.6666
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2016.6666
$ws.Range("I66").Value = 2016.6666
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10083.333
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6651.333000000001
$ws.Range("N66").ClearContents()
$ws.Range("H133").Value = 28795.084
$ws.Range("J133").Value = 28795.084
$ws.Range("L133").Value = 28795.084
$ws.Range("N133").Value = -33855.084
$ws.Range("H139").Value = 25980.883
$ws.Range("J139").Value = 25980.883
$ws.Range("L139").Value = 25980.883
$ws.Range("N139").Value = -36260.883

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 368.75
$ws.Range("I22").Value = 378.57144
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 378.57144
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -205.57144
$ws.Range("N22").Value = -646

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9187.25
$ws.Range("I86").Value = 8017.909
$ws.Range("J86").Value = 11759.8
$ws.Range("K86").Value = 8017.909
$ws.Range("L86").Value = 11759.8
$ws.Range("M86").Value = -6894.909
$ws.Range("N86").Value = -14005.8
$ws.Range("H89").Value = 9187.25
$ws.Range("I89").Value = 8017.909
$ws.Range("J89").Value = 11759.8
$ws.Range("K89").Value = 40089.545
$ws.Range("L89").Value = 58799
$ws.Range("M89").Value = -34473.545
$ws.Range("N89").Value = -70031

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 41506
$ws.Range("J35").Value = 41506
$ws.Range("L35").Value = 41506
$ws.Range("N35").Value = -42102
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H102").Value = 21851.34
$ws.Range("I102").Value = 2292.9697
$ws.Range("J102").Value = 54122.65
$ws.Range("K102").Value = 2292.9697
$ws.Range("L102").Value = 54122.65
$ws.Range("M102").Value = -670.9697000000001
$ws.Range("N102").Value = -57366.65
$ws.Range("H132").Value = 3427.158
$ws.Range("I132").Value = 3756.6
$ws.Range("K132").Value = 11269.8
$ws.Range("M132").Value = -8739.799999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 5674.778
$ws.Range("I31").Value = 7110.4287
$ws.Range("J31").Value = 650
$ws.Range("K31").Value = 7110.4287
$ws.Range("L31").Value = 650
$ws.Range("M31").Value = -6862.4287
$ws.Range("N31").Value = -1146
$ws.Range("H40").Value = 1688.1904
$ws.Range("I40").Value = 1197.1765
$ws.Range("J40").Value = 3775
$ws.Range("K40").Value = 1197.1765
$ws.Range("L40").Value = 3775
$ws.Range("M40").Value = -1061.1765
$ws.Range("N40").Value = -4047
$ws.Range("H46").Value = 1237.4642
$ws.Range("I46").Value = 920.7917
$ws.Range("K46").Value = 920.7917
$ws.Range("M46").Value = -732.7917
$ws.Range("H61").Value = 2937
$ws.Range("I61").Value = 1367.4445
$ws.Range("K61").Value = 1367.4445
$ws.Range("M61").Value = -1165.4445
$ws.Range("H94").Value = 25863.637
$ws.Range("J94").Value = 25863.637
$ws.Range("L94").Value = 25863.637
$ws.Range("N94").Value = -27215.637
$ws.Range("H113").Value = 2937
$ws.Range("I113").Value = 1367.4445
$ws.Range("K113").Value = 1367.4445
$ws.Range("M113").Value = 802.5554999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 56704
$ws.Range("J13").Value = 85006
$ws.Range("L13").Value = 85006
$ws.Range("N13").Value = -85286
$ws.Range("H69").Value = 37000
$ws.Range("J69").Value = 37000
$ws.Range("L69").Value = 37000
$ws.Range("N69").Value = -38498
$ws.Range("H72").Value = 37000
$ws.Range("J72").Value = 37000
$ws.Range("L72").Value = 111000
$ws.Range("N72").Value = -118488
$ws.Range("H88").Value = 44999
$ws.Range("J88").Value = 44999
$ws.Range("L88").Value = 44999
$ws.Range("N88").Value = -45811
$ws.Range("H91").Value = 44999
$ws.Range("J91").Value = 44999
$ws.Range("L91").Value = 44999
$ws.Range("N91").Value = -47807
$ws.Range("H126").Value = 2779748
$ws.Range("I126").Value = 1344.6086
$ws.Range("K126").Value = 4033.8258
$ws.Range("M126").Value = -1563.8258
$ws.Range("H136").Value = 1151.3158
$ws.Range("I136").Value = 541.6667
$ws.Range("K136").Value = 1625.0001
$ws.Range("M136").Value = 924.9999
